$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2023" column (S) by copying the formatting (and provisional
# values) of the existing 2022 column (R), then overwrite with the real
# 2023 figures.
$ws.Columns("R").Copy()
$ws.Columns("S").Insert()

$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 33.9
$ws.Range("S5").Value = 33.9
$ws.Range("S6").Value = 854
$ws.Range("S7").Value = 842
$ws.Range("S8").Value = 649.16999999999996
$ws.Range("S9").Value = 24.2
$ws.Range("S10").Value = 6.6
$ws.Range("S11").Value = 9.6999999999999993
$ws.Range("S12").Value = 0.8
$ws.Range("S13").Value = 24.1
$ws.Range("S14").Value = "_"

$ws.Range("D19").Select()
